# feat: add 2022-Q1 data
#
# - Insert a new worksheet "2022-Q1" (fund-holdings detail) right before
#   the "总计" (totals) summary sheet.
# - Insert a new top row into "总计" summarising the 2022-Q1 data.

$xlShiftDown    = [Microsoft.Office.Interop.Excel.XlInsertShiftDirection]::xlShiftDown
$xlPasteFormats = [Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats
$xlPasteValues  = [Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues

# Helper: force a numeric-looking string ("17.18", "011230", ...) to be
# stored as genuine text rather than being auto-coerced to a number
# (which would also silently drop leading zeros from fund codes).
function Set-TextValue {
    param($range, [string]$text)
    $escaped = $text.Replace('"', '""')
    $range.Formula = '="' + $escaped + '"'
    $range.Copy()
    $range.PasteSpecial($xlPasteValues)
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Add the "2022-Q1" sheet, positioned right before "总计".
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$newSheet = $wb.Worksheets.Add($totalSheet)
$newSheet.Name = "2022-Q1"

# Copy header-row formatting + text from an existing quarter sheet so the
# styling (bold, bordered, centered) matches exactly. (Column A of row 1
# is intentionally left blank, matching the template sheets.)
$templateSheet = $wb.Worksheets.Item("2021-Q3")
$templateSheet.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial($xlPasteFormats)

$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Row style for column A (bold, bordered) on the data rows, same as the
# other quarter sheets.
$templateSheet.Range("A2").Copy()
$newSheet.Range("A2:A8").PasteSpecial($xlPasteFormats)

$rows = @(
    @(0, "011230", "创金合信数字经济主题股票C", "17.18", "92.17", "5.52", "0.9483", 3),
    @(1, "011229", "创金合信数字经济主题股票A", "12.18", "92.17", "5.52", "0.6723", 3),
    @(2, "007497", "中庚价值灵动灵活配置混合", "24.35", "89.42", "2.20", "0.5357", 8),
    @(3, "501030", "汇添富中证环境治理指数（LOF）A", "6.61", "93.20", "2.29", "0.1514", 4),
    @(4, "501031", "汇添富中证环境治理指数（LOF）C", "2.74", "93.20", "2.29", "0.0627", 4),
    @(5, "164908", "交银施罗德中证环境治理指数（LOF）", "2.12", "93.72", "2.34", "0.0496", 3),
    @(6, "257050", "国联安主题驱动混合", "1.50", "65.37", "2.65", "0.0398", 10)
)

$r = 2
foreach ($row in $rows) {
    $newSheet.Range("A$r").Value = $row[0]
    Set-TextValue $newSheet.Range("B$r") $row[1]
    $newSheet.Range("C$r").Value = $row[2]
    Set-TextValue $newSheet.Range("D$r") $row[3]
    Set-TextValue $newSheet.Range("E$r") $row[4]
    Set-TextValue $newSheet.Range("F$r") $row[5]
    Set-TextValue $newSheet.Range("G$r") $row[6]
    $newSheet.Range("H$r").Value = $row[7]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# 2. Insert a new summary row at the top of "总计" for 2022-Q1.
# ---------------------------------------------------------------------
# NB: re-fetch by name — sheet references here track position, and
# "总计" moved from index 4 to index 5 once the new sheet was inserted
# before it.
$totalSheet = $wb.Worksheets.Item("总计")

$totalSheet.Range("A2:D2").Insert($xlShiftDown)
$totalSheet.Range("B2:D2").ClearFormats()

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 7
$totalSheet.Range("D2").Value = 2.46

# Restore column-A styling (bold, bordered) on the new row to match the
# other rows in the table (the Insert above leaves it unstyled).
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial($xlPasteFormats)

# Renumber the trailing index column (A) so it continues 0,1,2,3,...
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3

# Restore the original active sheet/selection (first sheet), since none
# of the above should change which tab the workbook opens to.
$wb.Worksheets.Item(1).Activate()
